$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "58.845.63"
$ws.Range("E2").Value = "  +2.54%  "

# Row 3
$ws.Range("D3").Value = "2.505.93"
$ws.Range("E3").Value = "  +3.32%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.32%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.51%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.09%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.568"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.02%  "

# Row 9
$ws.Range("D9").Value = "2.520.69"
$ws.Range("E9").Value = "  +3.57%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0995"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.53%  "

# Row 11
$ws.Range("E11").Value = "  -1.17%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.27"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.49%  "

# Row 13
$ws.Range("E13").Value = "  +2.19%  "

# Row 14
$ws.Range("D14").Value = "2.949.11"
$ws.Range("E14").Value = "  +3.16%  "

# Row 15
$ws.Range("D15").Value = "58.770.92"
$ws.Range("E15").Value = "  +2.51%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.43%  "

# Row 17
$ws.Range("E17").Value = "  +4.20%  "

# Row 18
$ws.Range("D18").Value = "2.521.20"
$ws.Range("E18").Value = "  +3.49%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.19%  "

# Row 20
$ws.Range("E20").Value = "  +4.88%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.13%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +10.55%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.12%  "

# Row 24
$ws.Range("E24").Value = "  +4.79%  "

# Row 25
$ws.Range("E25").Value = "  +2.83%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.53%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.161"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.56%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.37%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0763"
$ws.Range("E29").Value = "  +6.83%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.90%  "

# Row 31
$ws.Range("E31").Value = "  +5.82%  "

# Row 32
$ws.Range("E32").Value = "  +7.98%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.44%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.03%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.995"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.14%  "

# Row 36
$ws.Range("E36").Value = "  +3.74%  "

# Row 37
$ws.Range("E37").Value = "  +0.53%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.84%  "

# Row 39
$ws.Range("E39").Value = "  +5.97%  "

# Row 40
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.87%  "

# Row 41
$ws.Range("B41").Value = "SuiNetwork"
$ws.Range("C41").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.817"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.96%  "

# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.33%  "

# Row 43
$ws.Range("E43").Value = "  +5.15%  "

# Row 44
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "278.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.90%  "

# Row 45
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "132.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +11.99%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.595"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.48%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0936"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.37%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0512"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.58%  "

# Row 49
$ws.Range("E49").Value = "  +6.27%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.41%  "

# Row 51
$ws.Range("D51").Value = "1.755.95"
$ws.Range("E51").Value = "  +3.93%  "
